$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Ccl2"
$ws.Range("C2").Value = "Ccr3"
$ws.Range("D2").Value = "M1"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 4.025457
$ws.Range("H2").Value = 12.076371
$ws.Range("I2").Value = 0.007988230327331255
$ws.Range("J2").Value = 0.008053698920764607
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.1323866666666667
$ws.Range("N2").Value = 0.39716
$ws.Range("O2").Value = 0.2977240407890032
$ws.Range("P2").Value = 0.2977240407890032
$ws.Range("Q2").Value = 0.53291683404
$ws.Range("R2").Value = 4.796251506359999
$ws.Range("S2").Value = 0.002378288211806323
$ws.Range("T2").Value = 0.002397779785988072

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Ccl2"
$ws.Range("C3").Value = "Ccr3"
$ws.Range("D3").Value = "M2"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 4.025457
$ws.Range("H3").Value = 12.076371
$ws.Range("I3").Value = 0.007988230327331255
$ws.Range("J3").Value = 0.008053698920764607
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.3122756666666667
$ws.Range("N3").Value = 0.936827
$ws.Range("O3").Value = 0.7022759592109968
$ws.Range("P3").Value = 0.7022759592109967
$ws.Range("Q3").Value = 1.257052268313
$ws.Range("R3").Value = 11.313470414817
$ws.Range("S3").Value = 0.005609942115524932
$ws.Range("T3").Value = 0.005655919134776533

$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Ccl2"
$ws.Range("C4").Value = "Ccr3"
$ws.Range("D4").Value = "M1"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 143.2163033333333
$ws.Range("H4").Value = 429.64891
$ws.Range("I4").Value = 0.2842024688515132
$ws.Range("J4").Value = 0.2865316876050504
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.1323866666666667
$ws.Range("N4").Value = 0.39716
$ws.Range("O4").Value = 0.2977240407890032
$ws.Range("P4").Value = 0.2977240407890032
$ws.Range("Q4").Value = 18.95992901062222
$ws.Range("R4").Value = 170.6393610956
$ws.Range("S4").Value = 0.08461390742868331
$ws.Range("T4").Value = 0.08530737184786794

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Ccl2"
$ws.Range("C5").Value = "Ccr3"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 143.2163033333333
$ws.Range("H5").Value = 429.64891
$ws.Range("I5").Value = 0.2842024688515132
$ws.Range("J5").Value = 0.2865316876050504
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.3122756666666667
$ws.Range("N5").Value = 0.936827
$ws.Range("O5").Value = 0.7022759592109968
$ws.Range("P5").Value = 0.7022759592109967
$ws.Range("Q5").Value = 44.72296660095223
$ws.Range("R5").Value = 402.50669940857
$ws.Range("S5").Value = 0.1995885614228299
$ws.Range("T5").Value = 0.2012243157571824

$ws.Range("A6").Value = "M1"
$ws.Range("B6").Value = "Ccl2"
$ws.Range("C6").Value = "Ccr3"
$ws.Range("D6").Value = "M1"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 85.72041200000001
$ws.Range("H6").Value = 257.161236
$ws.Range("I6").Value = 0.1701060015818651
$ws.Range("J6").Value = 0.1715001277151632
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.1323866666666667
$ws.Range("N6").Value = 0.39716
$ws.Range("O6").Value = 0.2977240407890032
$ws.Range("P6").Value = 0.2977240407890032
$ws.Range("Q6").Value = 11.34823960997333
$ws.Range("R6").Value = 102.13415648976
$ws.Range("S6").Value = 0.05064464615341345
$ws.Range("T6").Value = 0.05105971101918848

$ws.Range("A7").Value = "M1"
$ws.Range("B7").Value = "Ccl2"
$ws.Range("C7").Value = "Ccr3"
$ws.Range("D7").Value = "M2"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 85.72041200000001
$ws.Range("H7").Value = 257.161236
$ws.Range("I7").Value = 0.1701060015818651
$ws.Range("J7").Value = 0.1715001277151632
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.3122756666666667
$ws.Range("N7").Value = 0.936827
$ws.Range("O7").Value = 0.7022759592109968
$ws.Range("P7").Value = 0.7022759592109967
$ws.Range("Q7").Value = 26.76839880424134
$ws.Range("R7").Value = 240.915589238172
$ws.Range("S7").Value = 0.1194613554284517
$ws.Range("T7").Value = 0.1204404166959746

$ws.Range("A8").Value = "M2"
$ws.Range("B8").Value = "Ccl2"
$ws.Range("C8").Value = "Ccr3"
$ws.Range("D8").Value = "M1"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 258.672133
$ws.Range("H8").Value = 776.016399
$ws.Range("I8").Value = 0.5133162713366615
$ws.Range("J8").Value = 0.5175232224251752
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.1323866666666667
$ws.Range("N8").Value = 0.39716
$ws.Range("O8").Value = 0.2977240407890032
$ws.Range("P8").Value = 0.2977240407890032
$ws.Range("Q8").Value = 34.24474144742666
$ws.Range("R8").Value = 308.2026730268399
$ws.Range("S8").Value = 0.1528265945050952
$ws.Range("T8").Value = 0.1540791049825692

$ws.Range("A9").Value = "M2"
$ws.Range("B9").Value = "Ccl2"
$ws.Range("C9").Value = "Ccr3"
$ws.Range("D9").Value = "M2"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 258.672133
$ws.Range("H9").Value = 776.016399
$ws.Range("I9").Value = 0.5133162713366615
$ws.Range("J9").Value = 0.5175232224251752
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.3122756666666667
$ws.Range("N9").Value = 0.936827
$ws.Range("O9").Value = 0.7022759592109968
$ws.Range("P9").Value = 0.7022759592109967
$ws.Range("Q9").Value = 80.77701278066365
$ws.Range("R9").Value = 726.993115025973
$ws.Range("S9").Value = 0.3604896768315662
$ws.Range("T9").Value = 0.363444117442606

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Ccl2"
$ws.Range("C10").Value = "Ccr3"
$ws.Range("D10").Value = "M1"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 12.2891965
$ws.Range("H10").Value = 24.578393
$ws.Range("I10").Value = 0.02438702790262897
$ws.Range("J10").Value = 0.0163912633338466
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.1323866666666667
$ws.Range("N10").Value = 0.39716
$ws.Range("O10").Value = 0.2977240407890032
$ws.Range("P10").Value = 0.2977240407890032
$ws.Range("Q10").Value = 1.626925760646666
$ws.Range("R10").Value = 9.761554563879999
$ws.Range("S10").Value = 0.007260604490004865
$ws.Range("T10").Value = 0.004880073153389436

$ws.Range("A11").Value = "sCs"
$ws.Range("B11").Value = "Ccl2"
$ws.Range("C11").Value = "Ccr3"
$ws.Range("D11").Value = "M2"
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 12.2891965
$ws.Range("H11").Value = 24.578393
$ws.Range("I11").Value = 0.02438702790262897
$ws.Range("J11").Value = 0.0163912633338466
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.3122756666666667
$ws.Range("N11").Value = 0.936827
$ws.Range("O11").Value = 0.7022759592109968
$ws.Range("P11").Value = 0.7022759592109967
$ws.Range("Q11").Value = 3.837617029835167
$ws.Range("R11").Value = 23.025702179011
$ws.Range("S11").Value = 0.0171264234126241
$ws.Range("T11").Value = 0.01151119018045716
